$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text columns (Coin name / Link URL) - safe to assign directly, Excel will not
# reinterpret these as numbers/dates.
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('B34').Value = 'LidoDAOToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('B50').Value = 'Aptos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'

# Numeric-looking / percent-looking text in Price & Volume(1h) columns.
# These must stay plain text (they already render as inline/shared strings in
# the source file, e.g. "1.000", "29.122.62", "  -1.37%  "), so force the cell
# to Text format before assigning, then restore the original style/format so no
# stray formatting is introduced.
$s = $ws.Range('D2').Style
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.122.62'
$ws.Range('D2').Style = $s
$s = $ws.Range('E2').Style
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -1.37%  '
$ws.Range('E2').Style = $s
$s = $ws.Range('D3').Style
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.836.48'
$ws.Range('D3').Style = $s
$s = $ws.Range('E3').Style
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -1.21%  '
$ws.Range('E3').Style = $s
$s = $ws.Range('D4').Style
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9991'
$ws.Range('D4').Style = $s
$s = $ws.Range('E4').Style
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('E4').Style = $s
$s = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.49'
$ws.Range('D5').Style = $s
$s = $ws.Range('E5').Style
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -2.13%  '
$ws.Range('E5').Style = $s
$s = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6650'
$ws.Range('D6').Style = $s
$s = $ws.Range('E6').Style
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -4.44%  '
$ws.Range('E6').Style = $s
$s = $ws.Range('E7').Style
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E7').Style = $s
$s = $ws.Range('D8').Style
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2957'
$ws.Range('D8').Style = $s
$s = $ws.Range('E8').Style
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -3.97%  '
$ws.Range('E8').Style = $s
$s = $ws.Range('D9').Style
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07360'
$ws.Range('D9').Style = $s
$s = $ws.Range('E9').Style
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -4.44%  '
$ws.Range('E9').Style = $s
$s = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '22.78'
$ws.Range('D10').Style = $s
$s = $ws.Range('E10').Style
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -3.77%  '
$ws.Range('E10').Style = $s
$s = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07680'
$ws.Range('D11').Style = $s
$s = $ws.Range('E11').Style
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -1.37%  '
$ws.Range('E11').Style = $s
$s = $ws.Range('D12').Style
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.844.51'
$ws.Range('D12').Style = $s
$s = $ws.Range('E12').Style
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.70%  '
$ws.Range('E12').Style = $s
$s = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6765'
$ws.Range('D14').Style = $s
$s = $ws.Range('E14').Style
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -2.84%  '
$ws.Range('E14').Style = $s
$s = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '86.31'
$ws.Range('D15').Style = $s
$s = $ws.Range('E15').Style
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -5.51%  '
$ws.Range('E15').Style = $s
$s = $ws.Range('D16').Style
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.223'
$ws.Range('D16').Style = $s
$s = $ws.Range('E16').Style
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -1.81%  '
$ws.Range('E16').Style = $s
$s = $ws.Range('D17').Style
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '29.050.80'
$ws.Range('D17').Style = $s
$s = $ws.Range('E17').Style
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -1.57%  '
$ws.Range('E17').Style = $s
$s = $ws.Range('D18').Style
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008231'
$ws.Range('D18').Style = $s
$s = $ws.Range('E18').Style
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.99%  '
$ws.Range('E18').Style = $s
$s = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '229.38'
$ws.Range('D19').Style = $s
$s = $ws.Range('E19').Style
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -3.71%  '
$ws.Range('E19').Style = $s
$s = $ws.Range('E20').Style
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -1.98%  '
$ws.Range('E20').Style = $s
$s = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.000'
$ws.Range('D21').Style = $s
$s = $ws.Range('E21').Style
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.03%  '
$ws.Range('E21').Style = $s
$s = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.312'
$ws.Range('D22').Style = $s
$s = $ws.Range('E22').Style
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -4.40%  '
$ws.Range('E22').Style = $s
$s = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9998'
$ws.Range('D23').Style = $s
$s = $ws.Range('E23').Style
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.02%  '
$ws.Range('E23').Style = $s
$s = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '161.09'
$ws.Range('D24').Style = $s
$s = $ws.Range('E24').Style
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.67%  '
$ws.Range('E24').Style = $s
$s = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1421'
$ws.Range('D25').Style = $s
$s = $ws.Range('E25').Style
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -4.97%  '
$ws.Range('E25').Style = $s
$s = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.693'
$ws.Range('D26').Style = $s
$s = $ws.Range('E26').Style
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -2.46%  '
$ws.Range('E26').Style = $s
$s = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.07'
$ws.Range('D27').Style = $s
$s = $ws.Range('E27').Style
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -1.21%  '
$ws.Range('E27').Style = $s
$s = $ws.Range('D28').Style
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.504'
$ws.Range('D28').Style = $s
$s = $ws.Range('E28').Style
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -2.02%  '
$ws.Range('E28').Style = $s
$s = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.231'
$ws.Range('D29').Style = $s
$s = $ws.Range('E30').Style
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -1.32%  '
$ws.Range('E30').Style = $s
$s = $ws.Range('E31').Style
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -1.65%  '
$ws.Range('E31').Style = $s
$s = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.05294'
$ws.Range('D32').Style = $s
$s = $ws.Range('E32').Style
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +3.41%  '
$ws.Range('E32').Style = $s
$s = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7508'
$ws.Range('D33').Style = $s
$s = $ws.Range('E33').Style
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -3.65%  '
$ws.Range('E33').Style = $s
$s = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.859'
$ws.Range('D34').Style = $s
$s = $ws.Range('E34').Style
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -1.44%  '
$ws.Range('E34').Style = $s
$s = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.132'
$ws.Range('D35').Style = $s
$s = $ws.Range('E35').Style
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -1.51%  '
$ws.Range('E35').Style = $s
$s = $ws.Range('E36').Style
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.15%  '
$ws.Range('E36').Style = $s
$s = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.315.62'
$ws.Range('D37').Style = $s
$s = $ws.Range('E37').Style
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.11%  '
$ws.Range('E37').Style = $s
$s = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01805'
$ws.Range('D38').Style = $s
$s = $ws.Range('E38').Style
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -3.77%  '
$ws.Range('E38').Style = $s
$s = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.715'
$ws.Range('D39').Style = $s
$s = $ws.Range('E39').Style
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.44%  '
$ws.Range('E39').Style = $s
$s = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9213'
$ws.Range('D40').Style = $s
$s = $ws.Range('E40').Style
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -2.91%  '
$ws.Range('E40').Style = $s
$s = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.997'
$ws.Range('D41').Style = $s
$s = $ws.Range('E41').Style
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +3.96%  '
$ws.Range('E41').Style = $s
$s = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9983'
$ws.Range('D42').Style = $s
$s = $ws.Range('E42').Style
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.29%  '
$ws.Range('E42').Style = $s
$s = $ws.Range('E43').Style
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -2.00%  '
$ws.Range('E43').Style = $s
$s = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.987.52'
$ws.Range('D44').Style = $s
$s = $ws.Range('E44').Style
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.57%  '
$ws.Range('E44').Style = $s
$s = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.5161'
$ws.Range('D45').Style = $s
$s = $ws.Range('E45').Style
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -1.29%  '
$ws.Range('E45').Style = $s
$s = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '63.89'
$ws.Range('D46').Style = $s
$s = $ws.Range('E46').Style
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +1.26%  '
$ws.Range('E46').Style = $s
$s = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.764'
$ws.Range('D47').Style = $s
$s = $ws.Range('E47').Style
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -1.36%  '
$ws.Range('E47').Style = $s
$s = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.286'
$ws.Range('D48').Style = $s
$s = $ws.Range('E48').Style
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -4.94%  '
$ws.Range('E48').Style = $s
$s = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05933'
$ws.Range('D49').Style = $s
$s = $ws.Range('E49').Style
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.01%  '
$ws.Range('E49').Style = $s
$s = $ws.Range('D50').Style
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.839'
$ws.Range('D50').Style = $s
$s = $ws.Range('E50').Style
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -1.87%  '
$ws.Range('E50').Style = $s
$s = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.07165'
$ws.Range('D51').Style = $s
$s = $ws.Range('E51').Style
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +6.53%  '
$ws.Range('E51').Style = $s
